# Insert a new data row for "Feria Lagunitas de Puerto Montt - Ajo" (weekly update).
# This shifts the existing row 470 (and everything below it) down by one row,
# and the new row 470 receives the latest week's price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 470 (pushes rows 470..511 to 471..512).
$ws.Rows.Item(470).Insert()

# Populate the newly inserted row 470 with the new weekly record.
$ws.Cells.Item(470, 1).Value = 4
$ws.Cells.Item(470, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(470, 3).Value = "Los Lagos"
$ws.Cells.Item(470, 4).Value = 45166
$ws.Cells.Item(470, 5).Value = 10
$ws.Cells.Item(470, 6).Value = 100112003
$ws.Cells.Item(470, 7).Value = "Ajo"
$ws.Cells.Item(470, 8).Value = "Chino"
$ws.Cells.Item(470, 9).Value = "Primera"
$ws.Cells.Item(470, 10).Value = 70
$ws.Cells.Item(470, 11).Value = 23000
$ws.Cells.Item(470, 12).Value = 23000
$ws.Cells.Item(470, 13).Value = 23000
$ws.Cells.Item(470, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(470, 15).Value = "China"
$ws.Cells.Item(470, 16).Value = 2300
$ws.Cells.Item(470, 17).Value = 10
$ws.Cells.Item(470, 18).Value = "Hortaliza"
